$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AK4").Value = 1040
